$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets: rename Sheet1, add two more ("LoginWithInvalidData",
# "HeartAttackCalculator") right after it, in order.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LoginWithValiddata"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LoginWithInvalidData"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "HeartAttackCalculator"

# ---------------------------------------------------------------------------
# Sheet 1: LoginWithValiddata — Username/Password header + one valid login
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "testuser@example.com"
$ws1.Range("B2").Value = 12345

$ws1.Range("A2:B2").Font.Size = 12
$ws1.Range("A2:B2").Font.Color = 0
$ws1.Rows.Item(2).RowHeight = 15.75

$ws1.Columns.Item(1).ColumnWidth = 22.5
$ws1.Columns.Item(2).ColumnWidth = 7.833333333333333

$ws1.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet 2: LoginWithInvalidData — header + one valid + two invalid (hyperlinked)
# login rows
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "testuser@example.com"
$ws2.Range("B2").Value = 123
$ws2.Range("A3").Value = "test@gmail.com"
$ws2.Range("B3").Value = 12345
$ws2.Range("A4").Value = "test@gmail.com"
$ws2.Range("B4").Value = 123456

$ws2.Range("A2:B2").Font.Size = 12
$ws2.Range("A2:B2").Font.Color = 0
$ws2.Rows.Item(2).RowHeight = 15.75

$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:test@gmail.com", "", "", "test@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:test@gmail.com", "", "", "test@gmail.com") | Out-Null

$ws2.Columns.Item(1).ColumnWidth = 22.5
$ws2.Columns.Item(2).ColumnWidth = 8.666666666666666

$ws2.Range("C4").Select()
$ws2.Activate()

# ---------------------------------------------------------------------------
# Sheet 3: HeartAttackCalculator — blank calculator sheet, formatting only
# ---------------------------------------------------------------------------
$ws3.Range("A2:B2").Font.Size = 12
$ws3.Range("A2:B2").Font.Color = 0
$ws3.Rows.Item(2).RowHeight = 15.75

$ws3.Columns.Item(1).ColumnWidth = 22.5

$ws3.PageSetup.Orientation = 1

$ws3.Range("A2").Select()

# Leave LoginWithInvalidData as the active sheet/tab, matching the target file.
$ws2.Activate()
